$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParaXml($doc, [string]$searchText, [string]$innerXml, [string]$pPrXml) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $searchText)
        return
    }
    $rng.Expand(4) | Out-Null
    $fullXml = "<w:p $wns>" + $pPrXml + $innerXml + "</w:p>"
    $rng.InsertXML($fullXml)
}

# ---------------------------------------------------------------
# 1. Names paragraph: split "Videet" and "Sek" out with spellStart/spellEnd
# ---------------------------------------------------------------
Replace-ParaXml $d "Honson Tran, Videet Parekh, Abelardo Lopez-Lagunas, Sek Chai" (
    '<w:r><w:t xml:space="preserve">Honson Tran, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Videet</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Parekh, Abelardo Lopez-Lagunas, </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Sek</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Chai</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 2. Table edits: "Raspberry Pi 4" -> "Raspberry Pi 4B (4GB)"
# ---------------------------------------------------------------
$tbl = $d.Tables.Item(1)

# Row 2 (ic_fp32): append a SECOND run "B (4GB)" rather than merging text
$row2Cell1 = $tbl.Rows.Item(2).Cells.Item(1)
$pPrTbl = (
    '<w:pPr>' +
    '<w:widowControl w:val="0"/>' +
    '<w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr>' +
    '<w:spacing w:line="240" w:lineRule="auto"/>' +
    '<w:jc w:val="center"/>' +
    '</w:pPr>'
)
$row2Xml = "<w:p $wns>" + $pPrTbl + '<w:r><w:t>Raspberry Pi 4</w:t></w:r><w:r><w:t>B (4GB)</w:t></w:r>' + "</w:p>"
$row2Cell1.Range.Paragraphs.Item(1).Range.InsertXML($row2Xml)

# Rows 3,4,5 (ic_int8, vww_fp32, vww_int8): simple in-place text replace (single run)
$tbl.Rows.Item(3).Cells.Item(1).Range.Text = "Raspberry Pi 4B (4GB)"
$tbl.Rows.Item(4).Cells.Item(1).Range.Text = "Raspberry Pi 4B (4GB)"
$tbl.Rows.Item(5).Cells.Item(1).Range.Text = "Raspberry Pi 4B (4GB)"

# ---------------------------------------------------------------
# 3. Add 4 new rows: kws_fp32, kws_int8, ad_fp32, ad_int8
# ---------------------------------------------------------------
$newRowsData = @(
    @("Raspberry Pi 4B (4GB)", "kws_fp32", "91.9% / 0.99"),
    @("Raspberry Pi 4B (4GB)", "kws_int8", "89.7% / 0.99"),
    @("Raspberry Pi 4B (4GB)", "ad_fp32", "83.1% / 0.90"),
    @("Raspberry Pi 4B (4GB)", "ad_int8", "77% / 0.86")
)

foreach ($rowData in $newRowsData) {
    $newRow = $tbl.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = $rowData[0]
    $newRow.Cells.Item(2).Range.Text = $rowData[1]
    $newRow.Cells.Item(3).Range.Text = $rowData[2]
}

Write-Output ("Table rows now: " + $tbl.Rows.Count)

# ---------------------------------------------------------------
# 4. lastRenderedPageBreak on "Yes (Visual Wake Words ... )"
# ---------------------------------------------------------------
$pPrList1 = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'
$ellipsis = [char]0x2026
$vwwSearch = "Yes (Visual Wake Words" + $ellipsis + " 80% Accuracy)"
Replace-ParaXml $d $vwwSearch (
    '<w:r><w:lastRenderedPageBreak/><w:t>Yes (Visual Wake Words' + $ellipsis + ' 80% Accuracy)</w:t></w:r>'
) $pPrList1

# ---------------------------------------------------------------
# 5. Keyword Spotting grammar split
# ---------------------------------------------------------------
$kwsSearch = "Yes (Keyword Spotting" + $ellipsis + " 90% Accuracy )"
Replace-ParaXml $d $kwsSearch (
    '<w:r><w:t xml:space="preserve">Yes (Keyword Spotting' + $ellipsis + ' 90% </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Accuracy )</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
) $pPrList1

# ---------------------------------------------------------------
# 6. "No, for some combination of benchmark, scenario and SUT" grammar split
# ---------------------------------------------------------------
Replace-ParaXml $d "No, for some combination of benchmark, scenario and SUT" (
    '<w:r><w:t xml:space="preserve">No, for some combination of benchmark, scenario and </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>SUT</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
) $pPrList1

# ---------------------------------------------------------------
# 7. "validation set in accuracy mode? (check one):"
# ---------------------------------------------------------------
Replace-ParaXml $d "For each SUT and benchmark, did the submission run on the whole validation set in accuracy mode? (check one):" (
    '<w:r><w:t xml:space="preserve">For each SUT and benchmark, did the submission run on the whole validation set in accuracy mode? (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> one):</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 8. "EEMBC Runner? (check one)"
# ---------------------------------------------------------------
Replace-ParaXml $d "For each SUT and benchmark, does the submission use the EEMBC Runner? (check one)" (
    '<w:r><w:t xml:space="preserve">For each SUT and benchmark, does the submission use the EEMBC Runner? (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> one)</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 9. "accuracy and performance modes? (check one)"
# ---------------------------------------------------------------
Replace-ParaXml $d "For each SUT and benchmark, is the same code run in accuracy and performance modes? (check one)" (
    '<w:r><w:t xml:space="preserve">For each SUT and benchmark, is the same code run in accuracy and performance modes? (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> one)</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 10. "outside of the official calibration set? (check one)"
# ---------------------------------------------------------------
Replace-ParaXml $d "Are the weights calibrated using data outside of the official calibration set? (check one)" (
    '<w:r><w:t xml:space="preserve">Are the weights calibrated using data outside of the official calibration set? (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> one)</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 11. "What numerics does the submission use? (check all that apply)"
# ---------------------------------------------------------------
Replace-ParaXml $d "What numerics does the submission use? (check all that apply)" (
    '<w:r><w:t xml:space="preserve">What </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>numerics</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> does the submission use? (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> all that apply)</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 12. "What backend does the submission use? (check all that apply)"
# ---------------------------------------------------------------
Replace-ParaXml $d "What backend does the submission use? (check all that apply)" (
    '<w:r><w:t xml:space="preserve">What backend does the submission use? (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> all that apply)</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 13. "Which of the following caching techniques ... (check all that apply, ideally none):"
# ---------------------------------------------------------------
Replace-ParaXml $d "Which of the following caching techniques does the submission use? (check all that apply, ideally none):" (
    '<w:r><w:t xml:space="preserve">Which of the following caching techniques does the submission use? (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> all that apply, ideally none):</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 14. "Which of the following techniques ... ideally none if submitting to the closed division.)"
# ---------------------------------------------------------------
Replace-ParaXml $d "Which of the following techniques does the submission use? (check all that apply, ideally none if submitting to the closed division.)" (
    '<w:r><w:t xml:space="preserve">Which of the following techniques does the submission use? (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> all that apply, ideally none if submitting to the closed division.)</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 15. "Is the submission congruent with all relevant MLPerf rules?"
# ---------------------------------------------------------------
Replace-ParaXml $d "Is the submission congruent with all relevant MLPerf rules?" (
    '<w:r><w:t xml:space="preserve">Is the submission congruent with all relevant </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>MLPerf</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> rules?</w:t></w:r>'
) ""

# ---------------------------------------------------------------
# 16. "Does your submission include the following: (check all that apply)"
# ---------------------------------------------------------------
Replace-ParaXml $d "Does your submission include the following: (check all that apply)" (
    '<w:r><w:t xml:space="preserve">Does your submission include the following: (check all that </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>apply)</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
) ""

# ---------------------------------------------------------------
# 17. "Scripts that set up and execute each system implementation tested"
# ---------------------------------------------------------------
$pPrList8 = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr>'
Replace-ParaXml $d "Scripts that set up and execute each system implementation tested" (
    '<w:r><w:t xml:space="preserve">Scripts that set up and execute each system implementation </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>tested</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
) $pPrList8

# ---------------------------------------------------------------
# 18. "Result logs for each system implementation tested"
# ---------------------------------------------------------------
Replace-ParaXml $d "Result logs for each system implementation tested" (
    '<w:r><w:t xml:space="preserve">Result logs for each system implementation </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>tested</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
) $pPrList8

Write-Output "All edits applied."
